# --- "7 Subat Testleri" sheet: new test-log entries added ---

$wb = $excel.ActiveWorkbook

# Excel alignment constants (xlHAlign.../xlVAlign...)
$xlLeft   = -4131
$xlCenter = -4108

# --- update selection on the previously-active sheet ("5 Subat Testleri") ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Range("C7").Select()

# --- add the new sheet after the last one and rename it ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "7 Subat Testleri"

# column widths
$ws.Columns.Item(2).ColumnWidth = 21.36
$ws.Columns.Item(3).ColumnWidth = 25.8

# --- row 1 ---
$ws.Range("A1").Value = "log0323"
$ws.Range("B1").Value = "1200-30 arasi git gelli 1 snlik periyotlu basincli test, dry run"
$ws.Range("C1").Value = "Motor istedigimiz periyodik hareketleri yapti, pozisyon ve hiz takibi beklendigi gibiydi. Dolu teste hazir"

$ws.Range("A1").VerticalAlignment = $xlCenter
$ws.Range("A1").HorizontalAlignment = $xlCenter

$ws.Range("B1").WrapText = $true
$ws.Range("B1").VerticalAlignment = $xlCenter
$ws.Range("B1").HorizontalAlignment = $xlCenter

$ws.Range("C1").WrapText = $true
$ws.Range("C1").VerticalAlignment = $xlCenter

$ws.Rows.Item(1).RowHeight = 58

# --- row 2 ---
$ws.Range("A2").Value = "log0325"
$ws.Range("B2").Value = "1200-30 arasi git gelli 1 snlik periyotlu basincli test"
$ws.Range("C2").Value = "Motor beklenen hareketi gerceklestirdi. 60 saniyenin sonunda 30 derece kapanma beklerken pos_feedback degeri 120 civari gosteriyordu. Ancak, biz elle kapatmayi deneyince tamamen kapaliydi, hizini alamamis olabilir mi?"

$ws.Range("A2").VerticalAlignment = $xlCenter

$ws.Range("B2").WrapText = $true
$ws.Range("B2").VerticalAlignment = $xlCenter
$ws.Range("B2").HorizontalAlignment = $xlCenter

$ws.Range("C2").WrapText = $true
$ws.Range("C2").VerticalAlignment = $xlCenter
$ws.Range("C2").HorizontalAlignment = $xlLeft

$ws.Rows.Item(2).RowHeight = 130.5

# --- row 3 ---
$ws.Range("A3").Value = "log0332"
$ws.Range("B3").Value = "600 derece ofsetli 80 derece magnitude 5Hz sinus fonksiyonudry run"
$ws.Range("C3").Value = "Test amacli dry run kosuldu"

$ws.Range("A3").VerticalAlignment = $xlCenter

$ws.Range("B3").WrapText = $true

$ws.Range("C3").WrapText = $true
$ws.Range("C3").VerticalAlignment = $xlCenter

$ws.Rows.Item(3).RowHeight = 43.5

# --- row 4 ---
$ws.Range("A4").Value = "log0335"
$ws.Range("B4").Value = "600 derece ofsetli 80 derece magnitude 5Hz sinus fonksiyonubasincli test"
$ws.Range("C4").Value = "Hiz 400RPMde 0.05sn kaldi tepe degeri olarak. Hiz ve pozisyon takibi iyi. Baglanti koptugundan log kaydetme ilginc bir sekilde yarida kesildi. Basinc sonda 20 bara kadar dustu. Test 17.5 sn surdu"

$ws.Range("A4").VerticalAlignment = $xlCenter

$ws.Range("B4").WrapText = $true
$ws.Range("B4").VerticalAlignment = $xlCenter

$ws.Range("C4").WrapText = $true
$ws.Range("C4").VerticalAlignment = $xlCenter

$ws.Rows.Item(4).RowHeight = 116

# --- row 5 ---
$ws.Range("A5").Value = "log0336"
$ws.Range("B5").Value = "600 derece ofsetli 40 derece magnitude 10Hz sinus fonksiyonubasincli test"
$ws.Range("C5").Value = "Dry Run"

$ws.Range("A5").VerticalAlignment = $xlCenter

$ws.Range("B5").WrapText = $true
$ws.Range("B5").VerticalAlignment = $xlCenter

$ws.Range("C5").WrapText = $true
$ws.Range("C5").VerticalAlignment = $xlCenter

$ws.Rows.Item(5).RowHeight = 58

# --- row 6 ---
$ws.Range("A6").Value = "log0339"
$ws.Range("B6").Value = "600 derece ofsetli 40 derece magnitude 10Hz sinus fonksiyonubasincli test"
$ws.Range("C6").Value = "Hiz 400RPMe ulasmadan azaliyor, istedigimiz gibi. Pozisyon ve hiz takibi iyi. Log basarili kaydoldu."

$ws.Range("A6").VerticalAlignment = $xlCenter

$ws.Range("B6").WrapText = $true
$ws.Range("B6").VerticalAlignment = $xlCenter

$ws.Range("C6").WrapText = $true
$ws.Range("C6").VerticalAlignment = $xlCenter

$ws.Rows.Item(6).RowHeight = 58

# --- row 7 ---
$ws.Range("A7").Value = "log0343"
$ws.Range("B7").Value = "600 derece ofsetli 25 derece magnitude 15Hz sinus fonksiyonu dry run"
$ws.Range("C7").Value = "Dry Run. Hiz demandi 400 olmasina ragmen motor hizi yakalayamadi. Bu testin ardindan 30 dereceli bir dry run kosuldu, onda hiz demandi 400RPM'de belli bir sure stabil kaldigindan dolayi 25 dereceye geri cektik. STLink baglanti kopmasaini onlemek icin USB soketinin ustune cift tarafli bant yapistirdik, titresimden dolayi baglanti kopuyordu."

$ws.Range("A7").VerticalAlignment = $xlCenter

$ws.Range("B7").WrapText = $true
$ws.Range("B7").VerticalAlignment = $xlCenter

$ws.Range("C7").WrapText = $true
$ws.Range("C7").VerticalAlignment = $xlCenter

$ws.Rows.Item(7).RowHeight = 188.5

# --- row 8 ---
$ws.Range("A8").Value = "log0345"
$ws.Range("B8").Value = "600 derece ofsetli 25 derece magnitude 15Hz sinus fonksiyonubasincli test"
$ws.Range("C8").Value = "Hiz demandini yakalayamadigi icin pozisyon takibi neredeyse 1 faz kadar gecikti. Sinus neredeyse cosinus'e donustu. Baglanti cift tarafli bant sayesinde kopmadi."

$ws.Range("A8").VerticalAlignment = $xlCenter

$ws.Range("B8").WrapText = $true
$ws.Range("B8").VerticalAlignment = $xlCenter

$ws.Range("C8").WrapText = $true
$ws.Range("C8").VerticalAlignment = $xlCenter

$ws.Rows.Item(8).RowHeight = 101.5

# --- selection / activation on the new sheet ---
$ws.Range("D5").Select()
